$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Cells.Item(2, 4)
$c.Value = "'35.428.57"
$c.Style = 'Normal'
$ws.Cells.Item(2, 5).Value = '  -3.88%  '

# Row 3
$c = $ws.Cells.Item(3, 4)
$c.Value = "'1.991.86"
$c.Style = 'Normal'
$ws.Cells.Item(3, 5).Value = '  -5.17%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.00%  '

# Row 5
$c = $ws.Cells.Item(5, 4)
$c.Value = "'240.69"
$c.Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  -1.99%  '

# Row 6
$c = $ws.Cells.Item(6, 4)
$c.Value = "'0.632"
$c.Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  -3.39%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.03%  '

# Row 8
$c = $ws.Cells.Item(8, 4)
$c.Value = "'56.16"
$c.Style = 'Normal'
$ws.Cells.Item(8, 5).Value = '  +3.00%  '

# Row 9
$c = $ws.Cells.Item(9, 4)
$c.Value = "'58.97"
$c.Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  -0.48%  '

# Row 10
$c = $ws.Cells.Item(10, 4)
$c.Value = "'0.355"
$c.Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  -3.55%  '

# Row 11
$c = $ws.Cells.Item(11, 4)
$c.Value = "'0.0720"
$c.Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  -6.12%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  -6.32%  '

# Row 13
$c = $ws.Cells.Item(13, 4)
$c.Value = "'0.894"
$c.Style = 'Normal'
$ws.Cells.Item(13, 5).Value = '  -5.99%  '

# Row 14
$c = $ws.Cells.Item(14, 4)
$c.Value = "'14.27"
$c.Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  -5.58%  '

# Row 15
$c = $ws.Cells.Item(15, 4)
$c.Value = "'2.273.48"
$c.Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  -5.49%  '

# Row 16
$c = $ws.Cells.Item(16, 4)
$c.Value = "'5.21"
$c.Style = 'Normal'
$ws.Cells.Item(16, 5).Value = '  -5.50%  '

# Row 17
$c = $ws.Cells.Item(17, 4)
$c.Value = "'1.981.81"
$c.Style = 'Normal'
$ws.Cells.Item(17, 5).Value = '  -4.41%  '

# Row 18
$c = $ws.Cells.Item(18, 4)
$c.Value = "'17.01"
$c.Style = 'Normal'
$ws.Cells.Item(18, 5).Value = '  -1.24%  '

# Row 19
$c = $ws.Cells.Item(19, 4)
$c.Value = "'35.430.54"
$c.Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  -3.81%  '

# Row 20
$c = $ws.Cells.Item(20, 4)
$c.Value = "'69.78"
$c.Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  -4.19%  '

# Row 21
$c = $ws.Cells.Item(21, 4)
$c.Value = "'0.0₃0833"
$c.Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  -5.61%  '

# Row 22
$c = $ws.Cells.Item(22, 4)
$c.Value = "'231.58"
$c.Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  -3.32%  '

# Row 23
$c = $ws.Cells.Item(23, 4)
$c.Value = "'5.01"
$c.Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  -8.75%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  +0.05%  '

# Row 25
$c = $ws.Cells.Item(25, 4)
$c.Value = "'2.26"
$c.Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  -5.85%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  +4.04%  '

# Row 27
$c = $ws.Cells.Item(27, 4)
$c.Value = "'9.12"
$c.Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  -5.69%  '

# Row 28
$c = $ws.Cells.Item(28, 4)
$c.Value = "'162.93"
$c.Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  -2.68%  '

# Row 29
$c = $ws.Cells.Item(29, 4)
$c.Value = "'19.43"
$c.Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  -7.86%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  -3.55%  '

# Row 31
$c = $ws.Cells.Item(31, 4)
$c.Value = "'1.14"
$c.Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  -2.23%  '

# Row 32
$c = $ws.Cells.Item(32, 4)
$c.Value = "'4.77"
$c.Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  -9.39%  '

# Row 33
$c = $ws.Cells.Item(33, 4)
$c.Value = "'0.0585"
$c.Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  -3.98%  '

# Row 34
$c = $ws.Cells.Item(34, 4)
$c.Value = "'0.0899"
$c.Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  +8.52%  '

# Row 35
$c = $ws.Cells.Item(35, 4)
$c.Value = "'4.24"
$c.Style = 'Normal'
$ws.Cells.Item(35, 5).Value = '  -10.51%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  -0.03%  '

# Row 37
$ws.Cells.Item(37, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Cells.Item(37, 4)
$c.Value = "'2.23"
$c.Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  -9.14%  '

# Row 38
$ws.Cells.Item(38, 2).Value = 'WEMIXToken'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$c = $ws.Cells.Item(38, 4)
$c.Value = "'1.80"
$c.Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  -2.51%  '

# Row 39
$c = $ws.Cells.Item(39, 4)
$c.Value = "'4.86"
$c.Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  -0.96%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  -7.68%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  -0.96%  '

# Row 42
$c = $ws.Cells.Item(42, 4)
$c.Value = "'0.0208"
$c.Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  -5.82%  '

# Row 43
$c = $ws.Cells.Item(43, 4)
$c.Value = "'1.08"
$c.Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  -7.14%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'Aave'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Cells.Item(45, 4)
$c.Value = "'90.13"
$c.Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  -6.91%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'Maker'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$c = $ws.Cells.Item(46, 4)
$c.Value = "'1.366.70"
$c.Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  -3.83%  '

# Row 47
$c = $ws.Cells.Item(47, 4)
$c.Value = "'7.41"
$c.Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  -5.25%  '

# Row 48
$c = $ws.Cells.Item(48, 4)
$c.Value = "'15.39"
$c.Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  -4.48%  '

# Row 49
$c = $ws.Cells.Item(49, 4)
$c.Value = "'2.89"
$c.Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  -0.73%  '

# Row 50
$c = $ws.Cells.Item(50, 4)
$c.Value = "'2.26"
$c.Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  -7.84%  '

# Row 51
$c = $ws.Cells.Item(51, 4)
$c.Value = "'45.14"
$c.Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  -2.02%  '
